$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1002539.2
$ws.Cells.Item(17, 10).Value = 1002539.2
$ws.Cells.Item(17, 12).Value = 3007617.6
$ws.Cells.Item(17, 14).Value = -3007953.6
$ws.Cells.Item(112, 8).Value = 1207.6
$ws.Cells.Item(112, 10).Value = 1176.8276
$ws.Cells.Item(112, 12).Value = 3530.4828
$ws.Cells.Item(112, 14).Value = -5746.4828
$ws.Cells.Item(137, 8).Value = 1321.9474
$ws.Cells.Item(137, 9).Value = 1352.55
$ws.Cells.Item(137, 10).Value = 1287.9445
$ws.Cells.Item(137, 11).Value = 4057.65
$ws.Cells.Item(137, 12).Value = 3863.8335
$ws.Cells.Item(137, 13).Value = -1507.65
$ws.Cells.Item(137, 14).Value = -8963.833500000001
$ws.Cells.Item(138, 8).Value = 3387.8591
$ws.Cells.Item(138, 10).Value = 4881.5674
$ws.Cells.Item(138, 12).Value = 14644.7022
$ws.Cells.Item(138, 14).Value = -24924.7022

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(16, 8).Value = 195
$ws.Cells.Item(16, 9).Value = 195
$ws.Cells.Item(16, 11).Value = 195
$ws.Cells.Item(16, 13).Value = 92
$ws.Cells.Item(32, 8).Value = 4750.805
$ws.Cells.Item(32, 9).Value = 3442.111
$ws.Cells.Item(32, 10).Value = 23596
$ws.Cells.Item(32, 11).Value = 3442.111
$ws.Cells.Item(32, 12).Value = 23596
$ws.Cells.Item(32, 13).Value = -3155.111
$ws.Cells.Item(32, 14).Value = -24170
$ws.Cells.Item(45, 8).Value = 1226.9
$ws.Cells.Item(45, 9).Value = 1075
$ws.Cells.Item(45, 10).Value = 1454.75
$ws.Cells.Item(45, 11).Value = 1075
$ws.Cells.Item(45, 12).Value = 1454.75
$ws.Cells.Item(45, 13).Value = -698
$ws.Cells.Item(45, 14).Value = -2208.75
$ws.Cells.Item(61, 8).Value = 1294.7142
$ws.Cells.Item(61, 9).Value = 806.3871
$ws.Cells.Item(61, 11).Value = 806.3871
$ws.Cells.Item(61, 13).Value = -594.3871
$ws.Cells.Item(74, 8).Value = 2948.2
$ws.Cells.Item(74, 9).Value = 3114.46
$ws.Cells.Item(74, 10).Value = 1285.6
$ws.Cells.Item(74, 11).Value = 3114.46
$ws.Cells.Item(74, 12).Value = 1285.6
$ws.Cells.Item(74, 13).Value = -2240.46
$ws.Cells.Item(74, 14).Value = -3033.6
$ws.Cells.Item(77, 8).Value = 2948.2
$ws.Cells.Item(77, 9).Value = 3114.46
$ws.Cells.Item(77, 10).Value = 1285.6
$ws.Cells.Item(77, 11).Value = 15572.3
$ws.Cells.Item(77, 12).Value = 6428
$ws.Cells.Item(77, 13).Value = -11204.3
$ws.Cells.Item(77, 14).Value = -15164
$ws.Cells.Item(102, 8).Value = 2266.6667
$ws.Cells.Item(102, 9).Value = 1900
$ws.Cells.Item(102, 11).Value = 1900
$ws.Cells.Item(102, 13).Value = -278
$ws.Cells.Item(136, 8).Value = 1294.7142
$ws.Cells.Item(136, 9).Value = 806.3871
$ws.Cells.Item(136, 11).Value = 2419.1613
$ws.Cells.Item(136, 13).Value = 130.8386999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1881.1177
$ws.Cells.Item(134, 9).Value = 1387.7273
$ws.Cells.Item(134, 11).Value = 4163.1819
$ws.Cells.Item(134, 13).Value = -1628.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2325.1018
$ws.Cells.Item(31, 9).Value = 1610.8788
$ws.Cells.Item(31, 10).Value = 3231.6155
$ws.Cells.Item(31, 11).Value = 1610.8788
$ws.Cells.Item(31, 12).Value = 3231.6155
$ws.Cells.Item(31, 13).Value = -1315.8788
$ws.Cells.Item(31, 14).Value = -3821.6155
$ws.Cells.Item(34, 8).Value = 2325.1018
$ws.Cells.Item(34, 9).Value = 1610.8788
$ws.Cells.Item(34, 10).Value = 3231.6155
$ws.Cells.Item(34, 11).Value = 1610.8788
$ws.Cells.Item(34, 12).Value = 3231.6155
$ws.Cells.Item(34, 13).Value = -1408.8788
$ws.Cells.Item(34, 14).Value = -3635.6155
$ws.Cells.Item(58, 8).Value = 1415.8727
$ws.Cells.Item(58, 9).Value = 1028.762
$ws.Cells.Item(58, 11).Value = 1028.762
$ws.Cells.Item(58, 13).Value = -825.7619999999999
$ws.Cells.Item(94, 8).Value = 3239.7144
$ws.Cells.Item(94, 9).Value = 3726
$ws.Cells.Item(94, 10).Value = 2875
$ws.Cells.Item(94, 11).Value = 3726
$ws.Cells.Item(94, 12).Value = 2875
$ws.Cells.Item(94, 13).Value = -3275
$ws.Cells.Item(94, 14).Value = -3777
$ws.Cells.Item(132, 8).Value = 1539.35
$ws.Cells.Item(132, 9).Value = 799.2121
$ws.Cells.Item(132, 10).Value = 5028.5713
$ws.Cells.Item(132, 11).Value = 2397.6363
$ws.Cells.Item(132, 12).Value = 15085.7139
$ws.Cells.Item(132, 13).Value = 132.3636999999999
$ws.Cells.Item(132, 14).Value = -20145.7139
$ws.Cells.Item(134, 8).Value = 1974.881
$ws.Cells.Item(134, 9).Value = 1762.6216
$ws.Cells.Item(134, 10).Value = 3545.6
$ws.Cells.Item(134, 11).Value = 5287.864799999999
$ws.Cells.Item(134, 12).Value = 10636.8
$ws.Cells.Item(134, 13).Value = -2752.864799999999
$ws.Cells.Item(134, 14).Value = -15706.8
$ws.Cells.Item(136, 8).Value = 1415.8727
$ws.Cells.Item(136, 9).Value = 1028.762
$ws.Cells.Item(136, 11).Value = 3086.286
$ws.Cells.Item(136, 13).Value = -536.2860000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 310.33334
$ws.Cells.Item(33, 9).Value = 244.25
$ws.Cells.Item(33, 10).Value = 343.375
$ws.Cells.Item(33, 11).Value = 1465.5
$ws.Cells.Item(33, 12).Value = 2060.25
$ws.Cells.Item(33, 13).Value = -1182.5
$ws.Cells.Item(33, 14).Value = -2626.25
$ws.Cells.Item(130, 8).Value = 1564.2142
$ws.Cells.Item(130, 9).Value = 1433
$ws.Cells.Item(130, 11).Value = 4299
$ws.Cells.Item(130, 13).Value = 721
$ws.Cells.Item(132, 8).Value = 1651.4667
$ws.Cells.Item(132, 9).Value = 1624.5714
$ws.Cells.Item(132, 10).Value = 1675
$ws.Cells.Item(132, 11).Value = 14621.1426
$ws.Cells.Item(132, 12).Value = 15075
$ws.Cells.Item(132, 13).Value = -12091.1426
$ws.Cells.Item(132, 14).Value = -20135
$ws.Cells.Item(140, 8).Value = 2156.3684
$ws.Cells.Item(140, 9).Value = 1798
$ws.Cells.Item(140, 10).Value = 2514.7368
$ws.Cells.Item(140, 11).Value = 5394
$ws.Cells.Item(140, 12).Value = 7544.2104
$ws.Cells.Item(140, 13).Value = -214
$ws.Cells.Item(140, 14).Value = -17904.2104

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1608
$ws.Cells.Item(113, 9).Value = 1445.6875
$ws.Cells.Item(113, 11).Value = 1445.6875
$ws.Cells.Item(113, 13).Value = 724.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 38134.266
$ws.Cells.Item(40, 9).Value = 45342.832
$ws.Cells.Item(40, 10).Value = 9300
$ws.Cells.Item(40, 11).Value = 45342.832
$ws.Cells.Item(40, 12).Value = 9300
$ws.Cells.Item(40, 13).Value = -45206.832
$ws.Cells.Item(40, 14).Value = -9572
$ws.Cells.Item(132, 8).Value = 5779.3286
$ws.Cells.Item(132, 9).Value = 5359.778
$ws.Cells.Item(132, 10).Value = 7195.3125
$ws.Cells.Item(132, 11).Value = 16079.334
$ws.Cells.Item(132, 12).Value = 21585.9375
$ws.Cells.Item(132, 13).Value = -13549.334
$ws.Cells.Item(132, 14).Value = -26645.9375
$ws.Cells.Item(136, 8).Value = 8773568
$ws.Cells.Item(136, 9).Value = 1693.0857
$ws.Cells.Item(136, 10).Value = 111112110
$ws.Cells.Item(136, 11).Value = 5079.257100000001
$ws.Cells.Item(136, 12).Value = 333336330
$ws.Cells.Item(136, 13).Value = -2529.257100000001
$ws.Cells.Item(136, 14).Value = -333341430

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1148.1666
$ws.Cells.Item(132, 9).Value = 980.9718
$ws.Cells.Item(132, 10).Value = 2061.3076
$ws.Cells.Item(132, 11).Value = 2942.9154
$ws.Cells.Item(132, 12).Value = 6183.9228
$ws.Cells.Item(132, 13).Value = -412.9153999999999
$ws.Cells.Item(132, 14).Value = -11243.9228
$ws.Cells.Item(136, 8).Value = 5377828.5
$ws.Cells.Item(136, 9).Value = 7246891
$ws.Cells.Item(136, 11).Value = 21740673
$ws.Cells.Item(136, 13).Value = -21738123
